$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 844.3077
$ws.Range("I28").Value = 865
$ws.Range("J28").Value = 797.75
$ws.Range("K28").Value = 865
$ws.Range("L28").Value = 797.75
$ws.Range("M28").Value = -380
$ws.Range("N28").Value = -1767.75
$ws.Range("H33").Value = 94.86667
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H113").Value = 9386.6875
$ws.Range("I113").Value = 12714.571
$ws.Range("J113").Value = 6798.3335
$ws.Range("K113").Value = 12714.571
$ws.Range("L113").Value = 6798.3335
$ws.Range("M113").Value = -9460.571
$ws.Range("N113").Value = -13306.3335
$ws.Range("H141").Value = 1462.25
$ws.Range("I141").Value = 749.6667
$ws.Range("J141").Value = 3600
$ws.Range("K141").Value = 2249.0001
$ws.Range("L141").Value = 10800
$ws.Range("M141").Value = 2930.9999
$ws.Range("N141").Value = -21160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1136.2222
$ws.Range("I132").Value = 1153.25
$ws.Range("K132").Value = 3459.75
$ws.Range("M132").Value = -929.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 24166.666
$ws.Range("J9").Value = 24166.666
$ws.Range("L9").Value = 24166.666
$ws.Range("N9").Value = -24502.666
$ws.Range("H99").Value = 1310
$ws.Range("I99").Value = 1328.2273
$ws.Range("K99").Value = 1328.2273
$ws.Range("M99").Value = 169.7727
$ws.Range("H134").Value = 1475
$ws.Range("I134").Value = 1475
$ws.Range("K134").Value = 4425
$ws.Range("M134").Value = -1890

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1116.3334
$ws.Range("I10").Value = 1116.3334
$ws.Range("K10").Value = 1116.3334
$ws.Range("M10").Value = -977.3334
$ws.Range("H41").Value = 11013
$ws.Range("I41").Value = 15000
$ws.Range("K41").Value = 15000
$ws.Range("M41").Value = -14572
$ws.Range("H51").Value = 15797.25
$ws.Range("I51").Value = 3090
$ws.Range("J51").Value = 20033
$ws.Range("K51").Value = 3090
$ws.Range("L51").Value = 20033
$ws.Range("M51").Value = -2354
$ws.Range("N51").Value = -21505
$ws.Range("H60").Value = 9000
$ws.Range("I60").Value = 3500
$ws.Range("K60").Value = 3500
$ws.Range("M60").Value = -2989
$ws.Range("H61").Value = 15797.25
$ws.Range("I61").Value = 3090
$ws.Range("J61").Value = 20033
$ws.Range("K61").Value = 3090
$ws.Range("L61").Value = 20033
$ws.Range("M61").Value = -2742
$ws.Range("N61").Value = -20729
$ws.Range("H86").Value = 11687326
$ws.Range("J86").Value = 5986
$ws.Range("L86").Value = 5986
$ws.Range("N86").Value = -8232
$ws.Range("H89").Value = 11687326
$ws.Range("J89").Value = 5986
$ws.Range("L89").Value = 29930
$ws.Range("N89").Value = -41162
$ws.Range("H93").Value = 4444
$ws.Range("I93").Value = 4444
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4444
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2572
$ws.Range("N93").ClearContents()
$ws.Range("H99").Value = 3967.875
$ws.Range("I99").Value = 2710.6
$ws.Range("J99").Value = 6063.3335
$ws.Range("K99").Value = 2710.6
$ws.Range("L99").Value = 6063.3335
$ws.Range("M99").Value = -1212.6
$ws.Range("N99").Value = -9059.333500000001
$ws.Range("H126").Value = 3967.875
$ws.Range("I126").Value = 2710.6
$ws.Range("J126").Value = 6063.3335
$ws.Range("K126").Value = 8131.799999999999
$ws.Range("L126").Value = 18190.0005
$ws.Range("M126").Value = -5661.799999999999
$ws.Range("N126").Value = -23130.0005
$ws.Range("H132").Value = 1528.8
$ws.Range("I132").Value = 1536.25
$ws.Range("K132").Value = 4608.75
$ws.Range("M132").Value = -2078.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 737.25
$ws.Range("I97").Value = 649.6667
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1949.0001
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1453.0001
$ws.Range("N97").Value = -3992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 10789.421
$ws.Range("I3").Value = 10625
$ws.Range("J3").Value = 11666.333
$ws.Range("K3").Value = 10625
$ws.Range("L3").Value = 11666.333
$ws.Range("M3").Value = -10509
$ws.Range("N3").Value = -11898.333
$ws.Range("H10").Value = 255998.25
$ws.Range("I10").Value = 255998.25
$ws.Range("K10").Value = 255998.25
$ws.Range("M10").Value = -255829.25
$ws.Range("H13").Value = 401
$ws.Range("I13").Value = 235
$ws.Range("J13").Value = 650
$ws.Range("K13").Value = 235
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = -96
$ws.Range("N13").Value = -928
$ws.Range("H113").Value = 4216.7144
$ws.Range("I113").Value = 3566.818
$ws.Range("K113").Value = 3566.818
$ws.Range("M113").Value = -1396.818
$ws.Range("H122").Value = 1381.75
$ws.Range("J122").Value = 1549.25
$ws.Range("L122").Value = 4647.75
$ws.Range("N122").Value = -9547.75
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1567.8889
$ws.Range("I7").Value = 1730.2
$ws.Range("J7").Value = 1365
$ws.Range("K7").Value = 1730.2
$ws.Range("L7").Value = 1365
$ws.Range("M7").Value = -1618.2
$ws.Range("N7").Value = -1589
$ws.Range("H22").Value = 1878
$ws.Range("I22").Value = 833.3333
$ws.Range("K22").Value = 833.3333
$ws.Range("M22").Value = -538.3333
$ws.Range("H27").Value = 1878
$ws.Range("I27").Value = 833.3333
$ws.Range("K27").Value = 833.3333
$ws.Range("M27").Value = -726.3333
$ws.Range("H82").Value = 1108.1538
$ws.Range("I82").Value = 1222.2222
$ws.Range("K82").Value = 1222.2222
$ws.Range("M82").Value = -861.2221999999999
$ws.Range("H85").Value = 1108.1538
$ws.Range("I85").Value = 1222.2222
$ws.Range("K85").Value = 1222.2222
$ws.Range("M85").Value = 25.77780000000007
$ws.Range("H100").Value = 4683.5
$ws.Range("I100").Value = 5614.6665
$ws.Range("K100").Value = 5614.6665
$ws.Range("M100").Value = -5073.6665
$ws.Range("H126").Value = 1567.8889
$ws.Range("I126").Value = 1730.2
$ws.Range("J126").Value = 1365
$ws.Range("K126").Value = 5190.6
$ws.Range("L126").Value = 4095
$ws.Range("M126").Value = -2720.6
$ws.Range("N126").Value = -9035
$ws.Range("H136").Value = 3333.3333
$ws.Range("I136").Value = 6000
$ws.Range("K136").Value = 18000
$ws.Range("M136").Value = -15450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 1533.5454
$ws.Range("I122").Value = 1613.4736
$ws.Range("K122").Value = 4840.4208
$ws.Range("M122").Value = -2390.4208
$ws.Range("H132").Value = 3138.889
$ws.Range("I132").Value = 2812.5
$ws.Range("K132").Value = 8437.5
$ws.Range("M132").Value = -5907.5
